$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.782259294303912

$ws.Range("C3").Value = -0.9921462019007898
$ws.Range("E3").Value = -0.3230348957779294

$ws.Range("C4").Value = 0.2467309912830284
$ws.Range("E4").Value = -0.06959526544320083

$ws.Range("C5").Value = 1.160201558804674
$ws.Range("E5").Value = 0.2932139896134167

$ws.Range("C6").Value = 1.048604932640185
$ws.Range("E6").Value = 1.078804187516891

$ws.Range("C7").Value = 1.06837811337479
$ws.Range("E7").Value = 1.25598608434605

$ws.Range("C8").Value = 1.435208340819005
$ws.Range("E8").Value = 1.407107513712802

$ws.Range("C9").Value = 1.577589817310243
$ws.Range("E9").Value = 1.464859320654099

$ws.Range("C10").Value = 1.979569114089963
$ws.Range("E10").Value = 1.639881111696151

$ws.Range("C11").Value = 1.815212363528707
$ws.Range("E11").Value = 1.806931013599544

$ws.Range("C12").Value = 2.284026378382942
$ws.Range("E12").Value = 2.119133965447961

$ws.Range("C13").Value = 1.264761787657309
$ws.Range("E13").Value = 1.798687504247187

$ws.Range("C14").Value = 0.824608016336259
$ws.Range("E14").Value = 1.395219579261608

$ws.Range("C15").Value = -1.4191429117966
$ws.Range("E15").Value = -0.1125839228000469

$ws.Range("C16").Value = 1.896944139870205
$ws.Range("E16").Value = 0.5116467003986713

$ws.Range("C17").Value = -0.7016063587211741
$ws.Range("E17").Value = 0.2942159770784825

$ws.Range("C18").Value = -0.1754728623905355
$ws.Range("E18").Value = -0.0331361487157622

$ws.Range("C19").Value = 0.5695821893874298
$ws.Range("E19").Value = 0.316149716722669
